# Apply updated dSF (column F) values per repull of data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -6
    3  = -4
    6  = -3
    7  = 7
    8  = 1
    10 = 4
    11 = 0
    13 = 1
    14 = 3
    16 = 3
    17 = 1
    19 = -1
    20 = -2
    21 = -2
    23 = -2
    25 = 7
    26 = -2
    27 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
